# Update "Pais" worksheet with refreshed COVID country/province stats.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1.
$ws.Range("A1").Value = "Datos actualizados a 19 de Octubre de 2020 a las 05:25"

# Row 26 - Pakistan
$ws.Range("B26").Value = 323452
$ws.Range("C26").Value = 433
$ws.Range("D26").Value = 307409
$ws.Range("E26").Value = 9384
$ws.Range("G26").Value = 5
$ws.Range("H26").Value = 6659

# Row 30 - Belgica
$ws.Range("B30").Value = 222253
$ws.Range("C30").Value = 9138
$ws.Range("D30").Value = 21157
$ws.Range("E30").Value = 190683
$ws.Range("G30").Value = 21
$ws.Range("H30").Value = 10413

# Row 45 - Kazajistan
$ws.Range("B45").Value = 109508
$ws.Range("C45").Value = 102
$ws.Range("E45").Value = 2739

# Row 153 - Belice
$ws.Range("B153").Value = 2813
$ws.Range("C153").Value = 38
$ws.Range("D153").Value = 1670
$ws.Range("E153").Value = 1099
$ws.Range("G153").Value = 1
$ws.Range("H153").Value = 44

# Row 174 - Islas Turcas y Caicos
$ws.Range("D174").Value = 684
$ws.Range("E174").Value = 8
